$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.936.53'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = '1.811.97'
$ws.Range("E3").Value = '  +2.02%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4290'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.21%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3692'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07229'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8620'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.14'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.39%  '

$ws.Range("D12").Value = '2.017.61'
$ws.Range("E12").Value = '  +8.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.626'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.08%  '

$ws.Range("E14").Value = '  +2.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06892'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.84%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '80.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.43%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008846'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.71%  '

$ws.Range("D21").Value = '26.983.23'
$ws.Range("E21").Value = '  -0.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.191'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.03%  '

$ws.Range("E23").Value = '  -1.16%  '

$ws.Range("D24").Value = '2.266.79'
$ws.Range("E24").Value = '  +10.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  -1.45%  '

$ws.Range("E27").Value = '  +0.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.222'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.897'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08904'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7424'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.69%  '

$ws.Range("E33").Value = '  +5.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.418'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.801'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.125'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.85%  '

$ws.Range("E38").Value = '  +2.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01920'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.72%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5079'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.740'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1644'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.425'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.248'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '106.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.58%  '

$ws.Range("E47").Value = '  -0.11%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4580'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.91%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.650'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06272'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.805'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.65%  '
